# "fixing matricula of Matc65"
# The "matricula" values (column A, text) for rows 14-39 were corrected.
# Apply the new values cell-by-cell, forcing text storage (the original
# cells are stored as text/inlineStr, not numbers) without leaving a
# residual number-format style on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newMatriculas = @{
    14 = "217216526"
    15 = "216117974"
    16 = "221117463"
    17 = "217125254"
    18 = "219218129"
    19 = "218215397"
    20 = "220117282"
    21 = "219217429"
    22 = "216216087"
    23 = "220121412"
    24 = "210201260"
    25 = "201520233"
    26 = "217117994"
    27 = "219118481"
    28 = "221119218"
    29 = "219215012"
    30 = "219121541"
    31 = "214007731"
    32 = "219215013"
    33 = "220117290"
    34 = "219118473"
    35 = "220117273"
    36 = "220120071"
    37 = "221216783"
    38 = "214120645"
    39 = "220217140"
}

foreach ($row in $newMatriculas.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force the numeric-looking string to be stored as text, matching the
    # original cell type, then drop the temporary number-format style so
    # no formatting difference is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $newMatriculas[$row]
    $cell.Style = "Normal"
}
